$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text in place before the column shift so it stays aligned with
# the "MODEL_CONDITION" header cell (currently column E, becomes column D).
$ws.Range("E1").Value = "MODELCONDITION"

# Remove the old column A (row index values 0/8/10) and shift B:F left into A:E.
$ws.Columns("A").Delete()
